# Remove the "reviews_count" header column (E1). Deleting the cell with a
# shift-left deletes E1 and shifts F1:K1 one column to the left, which both
# renumbers the remaining headers (F->E, G->F, ..., K->J) and shrinks the
# used range so the sheet dimension becomes A1:J1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlShiftToLeft = -4159
$ws.Cells.Item(1, 5).Delete($xlShiftToLeft)
